{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2026-02-21 Saturday\", \"2026-02-22 Sunday\"],\n  [\"495\u00f77=70, 5\", \"810\u00f79=90, 0\"],\n  [\"368\u00f78=46, 0\", \"575\u00f78=71, 7\"],\n  [\"775\u00f77=110, 5\", \"410\u00f79=45, 5\"],\n  [\"807\u00f77=115, 2\", \"158\u00f72=79, 0\"],\n  [\"781\u00f78=97, 5\", \"514\u00f77=73, 3\"],\n  [\"780\u00f78=97, 4\", \"378\u00f77=54, 0\"],\n  [\"673\u00f74=168, 1\", \"410\u00f78=51, 2\"],\n  [\"512\u00f79=56, 8\", \"274\u00f78=34, 2\"],\n  [\"675\u00f72=337, 1\", \"777\u00f76=129, 3\"],\n  [\"759\u00f72=379, 1\", \"156\u00f75=31, 1\"],\n  [\"211\u00f79=23, 4\", \"650\u00f76=108, 2\"],\n  [\"502\u00f78=62, 6\", \"997\u00f76=166, 1\"],\n  [\"821\u00f78=102, 5\", \"607\u00f75=121, 2\"],\n  [\"172\u00f72=86, 0\", \"319\u00f73=106, 1\"],\n  [\"784\u00f77=112, 0\", \"531\u00f79=59, 0\"],\n  [\"251\u00f75=50, 1\", \"377\u00f75=75, 2\"],\n  [\"811\u00f74=202, 3\", \"585\u00f79=65, 0\"],\n  [\"953\u00f79=105, 8\", \"313\u00f78=39, 1\"],\n  [\"623\u00f77=89, 0\", \"650\u00f77=92, 6\"],\n  [\"350\u00f77=50, 0\", \"830\u00f74=207, 2\"],\n  [\"812\u00f77=116, 0\", \"136\u00f74=34, 0\"],\n  [\"245\u00f74=61, 1\", \"246\u00f76=41, 0\"],\n  [\"703\u00f77=100, 3\", \"965\u00f79=107, 2\"],\n  [\"622\u00f75=124, 2\", \"923\u00f77=131, 6\"],\n  [\"480\u00f76=80, 0\", \"536\u00f76=89, 2\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + before);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n    @('2026-02-21 Saturday', '2026-02-22 Sunday'),\n    @('495\u00f77=70, 5', '810\u00f79=90, 0'),\n    @('368\u00f78=46, 0', '575\u00f78=71, 7'),\n    @('775\u00f77=110, 5', '410\u00f79=45, 5'),\n    @('807\u00f77=115, 2', '158\u00f72=79, 0'),\n    @('781\u00f78=97, 5', '514\u00f77=73, 3'),\n    @('780\u00f78=97, 4', '378\u00f77=54, 0'),\n    @('673\u00f74=168, 1', '410\u00f78=51, 2'),\n    @('512\u00f79=56, 8', '274\u00f78=34, 2'),\n    @('675\u00f72=337, 1', '777\u00f76=129, 3'),\n    @('759\u00f72=379, 1', '156\u00f75=31, 1'),\n    @('211\u00f79=23, 4', '650\u00f76=108, 2'),\n    @('502\u00f78=62, 6', '997\u00f76=166, 1'),\n    @('821\u00f78=102, 5', '607\u00f75=121, 2'),\n    @('172\u00f72=86, 0', '319\u00f73=106, 1'),\n    @('784\u00f77=112, 0', '531\u00f79=59, 0'),\n    @('251\u00f75=50, 1', '377\u00f75=75, 2'),\n    @('811\u00f74=202, 3', '585\u00f79=65, 0'),\n    @('953\u00f79=105, 8', '313\u00f78=39, 1'),\n    @('623\u00f77=89, 0', '650\u00f77=92, 6'),\n    @('350\u00f77=50, 0', '830\u00f74=207, 2'),\n    @('812\u00f77=116, 0', '136\u00f74=34, 0'),\n    @('245\u00f74=61, 1', '246\u00f76=41, 0'),\n    @('703\u00f77=100, 3', '965\u00f79=107, 2'),\n    @('622\u00f75=124, 2', '923\u00f77=131, 6'),\n    @('480\u00f76=80, 0', '536\u00f76=89, 2'),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}"}
